# Apply a cyclic rotation of per-observation data among specific row groups.
# Within each group, row N's data (columns A,B,E,F,G,H,Q,R,AC) is replaced by the
# data that used to be on the next row in the group (wrapping around at the end).
# Row-invariant columns (D,P,S,T,U,V,W,Y,AA,AD,AE,AG,AT,AW,AX,AY, etc.) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) whose values get rotated between rows in each group.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

# Row groups that participate in the rotation. Within a group, each row receives
# the values that originally belonged to the following row (cyclically).
$groups = @(
    , @(13, 14)
    , @(21, 22)
    , @(51, 52, 53, 54, 55, 56, 57)
)

foreach ($group in $groups) {
    # Snapshot the original values (and whether the cell even had content) for
    # every row in this group before mutating anything. Value2 is used because
    # the bare Value property getter is unreliable in this runtime.
    $snapshot = @{}
    foreach ($r in $group) {
        $rowData = @{}
        foreach ($col in $cols) {
            $cell = $ws.Range("$col$r")
            $rowData[$col] = $cell.Value2
        }
        $snapshot[$r] = $rowData
    }

    $count = $group.Count
    for ($i = 0; $i -lt $count; $i++) {
        $targetRow = $group[$i]
        $sourceRow = $group[($i + 1) % $count]
        $sourceData = $snapshot[$sourceRow]

        foreach ($col in $cols) {
            $ws.Range("$col$targetRow").Value2 = $sourceData[$col]
        }
    }
}
